# Reorder slides: move the slide currently at position 9 ("Efa ampy ve...")
# so that it becomes position 7, pushing the current slides 7 and 8 down by one.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$s.MoveTo(7)
